$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '25.018.40'
$ws.Range('E2').Value = '  -3.78%  '
Set-TextValue 'D3' '1.644.54'
$ws.Range('E3').Value = '  -5.60%  '
Set-TextValue 'D4' '0.9989'
Set-TextValue 'D5' '233.52'
$ws.Range('E5').Value = '  -5.55%  '
Set-TextValue 'D6' '1.000'
$ws.Range('E6').Value = '  +0.01%  '
Set-TextValue 'D7' '0.4776'
$ws.Range('E7').Value = '  -5.33%  '
Set-TextValue 'D8' '0.2602'
$ws.Range('E8').Value = '  -5.29%  '
Set-TextValue 'D9' '0.06107'
$ws.Range('E9').Value = '  -1.31%  '
Set-TextValue 'D10' '0.07038'
$ws.Range('E10').Value = '  -3.16%  '
Set-TextValue 'D11' '1.644.70'
$ws.Range('E11').Value = '  -5.99%  '
Set-TextValue 'D12' '14.65'
$ws.Range('E12').Value = '  -3.13%  '
Set-TextValue 'D13' '0.5959'
$ws.Range('E13').Value = '  -9.02%  '
$ws.Range('E14').Value = '  -6.77%  '
Set-TextValue 'D15' '73.67'
$ws.Range('E15').Value = '  -5.06%  '
Set-TextValue 'D16' '0.9999'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('E17').Value = '  +0.03%  '
Set-TextValue 'D18' '25.021.56'
$ws.Range('E18').Value = '  -3.83%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.000006598'
$ws.Range('E19').Value = '  -3.72%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D20' '11.27'
$ws.Range('E20').Value = '  -5.55%  '
Set-TextValue 'D21' '1.855.87'
$ws.Range('E21').Value = '  -6.07%  '
Set-TextValue 'D22' '4.329'
$ws.Range('E22').Value = '  -3.24%  '
Set-TextValue 'D23' '8.577'
$ws.Range('E23').Value = '  -1.82%  '
Set-TextValue 'D24' '5.247'
$ws.Range('E24').Value = '  -2.82%  '
Set-TextValue 'D25' '134.46'
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('E26').Value = '  -2.32%  '
Set-TextValue 'D27' '1.388'
$ws.Range('E27').Value = '  -7.80%  '
Set-TextValue 'D28' '103.70'
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('E29').Value = '  -8.14%  '
Set-TextValue 'D30' '3.879'
$ws.Range('E30').Value = '  -1.71%  '
Set-TextValue 'D31' '0.07695'
$ws.Range('E31').Value = '  -5.84%  '
Set-TextValue 'D32' '3.568'
$ws.Range('E32').Value = '  -3.16%  '
Set-TextValue 'D33' '0.9990'
$ws.Range('E33').Value = '  +0.01%  '
Set-TextValue 'D34' '0.04293'
$ws.Range('E34').Value = '  -8.41%  '
Set-TextValue 'D35' '2.571'
$ws.Range('E35').Value = '  -3.21%  '
Set-TextValue 'D36' '0.5940'
$ws.Range('E36').Value = '  -3.12%  '
Set-TextValue 'D37' '0.9286'
$ws.Range('E37').Value = '  -6.90%  '
$ws.Range('E38').Value = '  -6.39%  '
Set-TextValue 'D39' '0.8778'
$ws.Range('E39').Value = '  +15.08%  '
$ws.Range('E40').Value = '  -0.03%  '
Set-TextValue 'D41' '0.01511'
$ws.Range('E41').Value = '  -6.95%  '
Set-TextValue 'D42' '98.86'
$ws.Range('E42').Value = '  -2.01%  '
Set-TextValue 'D43' '1.767'
$ws.Range('E43').Value = '  -8.34%  '
Set-TextValue 'D44' '0.3707'
$ws.Range('E44').Value = '  -5.46%  '
Set-TextValue 'D45' '4.675'
$ws.Range('E45').Value = '  -6.81%  '
$ws.Range('E46').Value = '  -5.34%  '
Set-TextValue 'D47' '6.109'
$ws.Range('E47').Value = '  -3.25%  '
Set-TextValue 'D48' '0.05211'
$ws.Range('E48').Value = '  -1.66%  '
Set-TextValue 'D49' '29.04'
$ws.Range('E49').Value = '  -5.49%  '
Set-TextValue 'D50' '0.9996'
$ws.Range('E50').Value = '  -0.23%  '
Set-TextValue 'D51' '0.9970'
$ws.Range('E51').Value = '  -0.28%  '
